$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the sub-header labels in row 2: the "total" columns (B and F) were
# mislabeled with leftover pandas placeholder text.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"

# Remove the two placeholder rows that carried a category label but no data
# ("situação do domicílio" and "grandes regiões e unidades da federação"),
# letting the real data rows shift up to close the gaps.
$ws.Rows.Item(8).EntireRow.Delete()
$ws.Rows.Item(5).EntireRow.Delete()
